$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.733.60"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3
$ws.Range("D3").Value = "2.221.70"
$ws.Range("E3").Value = "  -5.43%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'295.78"
$ws.Range("E5").Value = "  -4.89%  "

# Row 6
$ws.Range("D6").Value = "'84.09"
$ws.Range("E6").Value = "  -2.00%  "

# Row 7
$ws.Range("E7").Value = "  -2.80%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("E9").Value = "  -3.74%  "

# Row 10
$ws.Range("D10").Value = "'0.0785"
$ws.Range("E10").Value = "  -3.43%  "

# Row 11
$ws.Range("D11").Value = "'29.79"
$ws.Range("E11").Value = "  -1.21%  "

# Row 12
$ws.Range("D12").Value = "'47.48"
$ws.Range("E12").Value = "  -9.52%  "

# Row 13
$ws.Range("E13").Value = "  -2.25%  "

# Row 14
$ws.Range("D14").Value = "2.562.73"
$ws.Range("E14").Value = "  -5.44%  "

# Row 15
$ws.Range("D15").Value = "'6.30"
$ws.Range("E15").Value = "  -2.12%  "

# Row 16
$ws.Range("D16").Value = "'14.16"
$ws.Range("E16").Value = "  -4.51%  "

# Row 17
$ws.Range("D17").Value = "2.222.34"
$ws.Range("E17").Value = "  -5.70%  "

# Row 18
$ws.Range("D18").Value = "'0.721"
$ws.Range("E18").Value = "  -5.27%  "

# Row 19
$ws.Range("D19").Value = "39.616.71"
$ws.Range("E19").Value = "  -1.29%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0883"
$ws.Range("E20").Value = "  -2.43%  "

# Row 21
$ws.Range("D21").Value = "'5.76"
$ws.Range("E21").Value = "  -5.84%  "

# Row 22
$ws.Range("D22").Value = "'65.12"
$ws.Range("E22").Value = "  -4.60%  "

# Row 23
$ws.Range("D23").Value = "'10.49"
$ws.Range("E23").Value = "  -2.56%  "

# Row 24
$ws.Range("D24").Value = "'232.59"
$ws.Range("E24").Value = "  -1.33%  "

# Row 25
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("E26").Value = "  -5.38%  "

# Row 27
$ws.Range("E27").Value = "  -0.38%  "

# Row 28
$ws.Range("D28").Value = "'22.85"
$ws.Range("E28").Value = "  -3.19%  "

# Row 29
$ws.Range("E29").Value = "  +2.44%  "

# Row 30
$ws.Range("E30").Value = "  -1.07%  "

# Row 31
$ws.Range("D31").Value = "'32.33"
$ws.Range("E31").Value = "  -6.88%  "

# Row 32
$ws.Range("D32").Value = "'149.71"
$ws.Range("E32").Value = "  -2.66%  "

# Row 33
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.24%  "

# Row 34
$ws.Range("E34").Value = "  -5.65%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0706"
$ws.Range("E35").Value = "  -2.34%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  -3.00%  "

# Row 37
$ws.Range("D37").Value = "'16.07"
$ws.Range("E37").Value = "  +2.63%  "

# Row 38
$ws.Range("E38").Value = "  -2.22%  "

# Row 39
$ws.Range("E39").Value = "  -1.42%  "

# Row 40
$ws.Range("D40").Value = "'2.67"
$ws.Range("E40").Value = "  -5.95%  "

# Row 41
$ws.Range("D41").Value = "'1.66"
$ws.Range("E41").Value = "  -4.13%  "

# Row 42
$ws.Range("D42").Value = "'3.69"
$ws.Range("E42").Value = "  -5.10%  "

# Row 43
$ws.Range("D43").Value = "1.940.44"
$ws.Range("E43").Value = "  -1.12%  "

# Row 44
$ws.Range("E44").Value = "  -3.84%  "

# Row 45
$ws.Range("D45").Value = "'0.0267"
$ws.Range("E45").Value = "  +0.30%  "

# Row 46
$ws.Range("D46").Value = "'9.41"
$ws.Range("E46").Value = "  +0.30%  "

# Row 47
$ws.Range("D47").Value = "'16.34"
$ws.Range("E47").Value = "  -7.83%  "

# Row 48
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -4.21%  "

# Row 49
$ws.Range("D49").Value = "2.430.61"
$ws.Range("E49").Value = "  -5.40%  "

# Row 50
$ws.Range("D50").Value = "'71.11"
$ws.Range("E50").Value = "  +0.52%  "

# Row 51
$ws.Range("D51").Value = "'89.01"
$ws.Range("E51").Value = "  -4.66%  "
